$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price/volume table with the day's refreshed figures.
# NOTE: "Price" values that are plain decimal numbers (e.g. "581.20") would
# otherwise be auto-converted to numeric cells by Excel's type inference, but
# the source data stores them as literal text. We force text by prefixing
# with an apostrophe (Excel's standard "quote-prefix" mechanism) and then
# reset the cell style back to Normal so no stray formatting is left behind.
$ws.Range("D2").Value = "66.918.89"
$ws.Range("E2").Value = "  -1.66%  "

$ws.Range("D3").Value = "2.466.11"
$ws.Range("E3").Value = "  -2.94%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'581.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.85%  "

$ws.Range("D6").Value = "'168.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.25%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").Value = "'0.512"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.40%  "

$ws.Range("D9").Value = "2.467.34"
$ws.Range("E9").Value = "  -2.70%  "

$ws.Range("E10").Value = "  -3.09%  "

$ws.Range("E11").Value = "  -0.96%  "

$ws.Range("D12").Value = "'4.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.59%  "

$ws.Range("D13").Value = "'0.328"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.78%  "

$ws.Range("D14").Value = "2.916.72"
$ws.Range("E14").Value = "  -3.35%  "

$ws.Range("D15").Value = "'25.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.13%  "

$ws.Range("D16").Value = "66.870.45"
$ws.Range("E16").Value = "  -1.52%  "

$ws.Range("D17").Value = "'0.0000168"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.91%  "

$ws.Range("D18").Value = "2.470.94"
$ws.Range("E18").Value = "  -2.07%  "

$ws.Range("D19").Value = "'10.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.18%  "

$ws.Range("D20").Value = "'7.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.15%  "

$ws.Range("D21").Value = "'349.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.92%  "

$ws.Range("D22").Value = "'4.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.41%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").Value = "'68.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.61%  "

$ws.Range("D25").Value = "'4.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.86%  "

$ws.Range("D26").Value = "'1.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.61%  "

$ws.Range("D27").Value = "'9.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.91%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -58.72%  "

$ws.Range("D30").Value = "0.0₃0897"
$ws.Range("E30").Value = "  -7.43%  "

$ws.Range("D31").Value = "'505.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.69%  "

$ws.Range("D32").Value = "'7.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.29%  "

$ws.Range("E33").Value = "  -5.88%  "

$ws.Range("D34").Value = "'1.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.64%  "

$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("D36").Value = "'159.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.37%  "

$ws.Range("D37").Value = "'0.114"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.36%  "

$ws.Range("E38").Value = "  +0.21%  "

$ws.Range("D39").Value = "'18.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.52%  "

$ws.Range("D40").Value = "'1.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.54%  "

$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.88%  "

$ws.Range("D43").Value = "'4.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.21%  "

$ws.Range("D44").Value = "'0.325"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.12%  "

$ws.Range("D45").Value = "'2.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.27%  "

$ws.Range("D46").Value = "'38.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.17%  "

$ws.Range("D47").Value = "'140.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.26%  "

$ws.Range("D48").Value = "'3.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.89%  "

$ws.Range("D49").Value = "'0.509"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.73%  "

$ws.Range("D50").Value = "0.0₆0252"
$ws.Range("E50").Value = "  -11.51%  "

$ws.Range("E51").Value = "  -8.28%  "
